$d = $word.ActiveDocument

# --- Change 1: split "têm dados como uma data, um " run into three runs ---
$d.Content.Find.Execute("têm dados como uma data, um ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "têm dados como uma condição, uma data, um ", 2)

Write-Output "done"
